$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.00806623659493
$ws.Range("D2").Value = 1.027906497288259
$ws.Range("E2").Value = 1.010674910026534
$ws.Range("F2").Value = 1.022979768844466
$ws.Range("I2").Value = 1.029697109906932
$ws.Range("J2").Value = 1.013334225431748
$ws.Range("K2").Value = 1.030724861401493
$ws.Range("L2").Value = 1.013544385695603
$ws.Range("M2").Value = 1.025812555103948
$ws.Range("N2").Value = 1.008573996590691

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.009046869467235
$ws.Range("D3").Value = 1.028273613981299
$ws.Range("E3").Value = 1.011507227475956
$ws.Range("F3").Value = 1.024145426759398
$ws.Range("I3").Value = 1.029734777232134
$ws.Range("J3").Value = 1.013946452146486
$ws.Range("K3").Value = 1.030900992768428
$ws.Range("L3").Value = 1.014180750738654
$ws.Range("M3").Value = 1.026784015325728
$ws.Range("N3").Value = 1.00878141055065

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.009681980273114
$ws.Range("D4").Value = 1.028510937619181
$ws.Range("E4").Value = 1.012046661634972
$ws.Range("F4").Value = 1.024899744210854
$ws.Range("I4").Value = 1.02975759793294
$ws.Range("J4").Value = 1.014342636655408
$ws.Range("K4").Value = 1.031013990097006
$ws.Range("L4").Value = 1.014592750354448
$ws.Range("M4").Value = 1.027412108279556
$ws.Range("N4").Value = 1.008915492840075

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.009949118429452
$ws.Range("D5").Value = 1.028610652509268
$ws.Range("E5").Value = 1.012273646972577
$ws.Range("F5").Value = 1.025216873300158
$ws.Range("I5").Value = 1.029766819342045
$ws.Range("J5").Value = 1.014509200048662
$ws.Range("K5").Value = 1.031061260318145
$ws.Range("L5").Value = 1.014766009207686
$ws.Range("M5").Value = 1.027676036966348
$ws.Range("N5").Value = 1.008971830032323

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.009993980161153
$ws.Range("D6").Value = 1.028627391737629
$ws.Range("E6").Value = 1.012311770936157
$ws.Range("F6").Value = 1.025270121491079
$ws.Range("I6").Value = 1.029768345797604
$ws.Range("J6").Value = 1.014537167197786
$ws.Range("K6").Value = 1.031069183445186
$ws.Range("L6").Value = 1.014795103293874
$ws.Range("M6").Value = 1.027720344616224
$ws.Range("N6").Value = 1.008981287469251

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.009685549245985
$ws.Range("D7").Value = 1.028512270236885
$ws.Range("E7").Value = 1.012049693811009
$ws.Range("F7").Value = 1.024903981650223
$ws.Range("I7").Value = 1.029757722614105
$ws.Range("J7").Value = 1.014344862255004
$ws.Range("K7").Value = 1.031014622644546
$ws.Range("L7").Value = 1.014595065235426
$ws.Range("M7").Value = 1.027415635385727
$ws.Range("N7").Value = 1.008916245742539

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.008397526688504
$ws.Range("D8").Value = 1.028030611565365
$ws.Range("E8").Value = 1.01095601513833
$ws.Range("F8").Value = 1.023373696180551
$ws.Range("I8").Value = 1.029710160793349
$ws.Range("J8").Value = 1.013541123237926
$ws.Range("K8").Value = 1.030784586191858
$ws.Range("L8").Value = 1.013759400501455
$ws.Range("M8").Value = 1.026140969402435
$ws.Range("N8").Value = 1.00864411954181

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.006132287613383
$ws.Range("D9").Value = 1.027180229211031
$ws.Range("E9").Value = 1.009035507368092
$ws.Range("F9").Value = 1.020677590000688
$ws.Range("I9").Value = 1.029614491710792
$ws.Range("J9").Value = 1.012125103867247
$ws.Range("K9").Value = 1.030371849220165
$ws.Range("L9").Value = 1.012288630880148
$ws.Range("M9").Value = 1.023890975209704
$ws.Range("N9").Value = 1.00816362634608

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.004625114266183
$ws.Range("D10").Value = 1.026612325174861
$ws.Range("E10").Value = 1.007759718585061
$ws.Range("F10").Value = 1.018880488782785
$ws.Range("I10").Value = 1.02954278166494
$ws.Range("J10").Value = 1.011181298605331
$ws.Range("K10").Value = 1.03009180356355
$ws.Range("L10").Value = 1.01130934665513
$ws.Range("M10").Value = 1.022388389553086
$ws.Range("N10").Value = 1.007842658972998

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.00397320194103
$ws.Range("D11").Value = 1.026366208599601
$ws.Range("E11").Value = 1.007208377073682
$ws.Range("F11").Value = 1.018102394537615
$ws.Range("I11").Value = 1.029509858319811
$ws.Range("J11").Value = 1.010772674633538
$ws.Range("K11").Value = 1.029969397534717
$ws.Range("L11").Value = 1.010885603851498
$ws.Range("M11").Value = 1.021737140797976
$ws.Range("N11").Value = 1.007703528117989

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.00373115848209
$ws.Range("D12").Value = 1.026274760382811
$ws.Range("E12").Value = 1.007003747827305
$ws.Range("F12").Value = 1.017813384815855
$ws.Range("I12").Value = 1.029497348511827
$ws.Range("J12").Value = 1.010620901337706
$ws.Range("K12").Value = 1.029923759845133
$ws.Range("L12").Value = 1.01072825152292
$ws.Range("M12").Value = 1.02149514498474
$ws.Range("N12").Value = 1.007651826433479

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.003783072836349
$ws.Range("D13").Value = 1.026294377640812
$ws.Range("E13").Value = 1.007047634082525
$ws.Range("F13").Value = 1.01787537797338
$ws.Range("I13").Value = 1.029500044595197
$ws.Range("J13").Value = 1.010653456864576
$ws.Range("K13").Value = 1.029933556985607
$ws.Range("L13").Value = 1.010762002104766
$ws.Range("M13").Value = 1.021547058117788
$ws.Range("N13").Value = 1.007662917625374

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.003953192394305
$ws.Range("D14").Value = 1.026358650052763
$ws.Range("E14").Value = 1.007191459015773
$ws.Range("F14").Value = 1.018078504719175
$ws.Range("I14").Value = 1.029508829975547
$ws.Range("J14").Value = 1.010760128843632
$ws.Range("K14").Value = 1.029965628582615
$ws.Range("L14").Value = 1.010872596145965
$ws.Range("M14").Value = 1.021717139246284
$ws.Range("N14").Value = 1.00769925489403

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.004058022691451
$ws.Range("D15").Value = 1.026398246550641
$ws.Range("E15").Value = 1.007280095997364
$ws.Range("F15").Value = 1.018203659012882
$ws.Range("I15").Value = 1.029514205775396
$ws.Range("J15").Value = 1.010825854019764
$ws.Range("K15").Value = 1.029985366379101
$ws.Range("L15").Value = 1.010940742711412
$ws.Range("M15").Value = 1.021821919515923
$ws.Range("N15").Value = 1.007721640537603

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.004668394373045
$ws.Range("D16").Value = 1.02662865486379
$ws.Range("E16").Value = 1.007796332197576
$ws.Range("F16").Value = 1.018932129628978
$ws.Range("I16").Value = 1.029544927288855
$ws.Range("J16").Value = 1.011208418713553
$ws.Range("K16").Value = 1.030099903247808
$ws.Range("L16").Value = 1.011337475299468
$ws.Range("M16").Value = 1.022431597741799
$ws.Range("N16").Value = 1.007851889491391

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.005051452856726
$ws.Range("D17").Value = 1.026773129034443
$ws.Range("E17").Value = 1.008120444136522
$ws.Range("F17").Value = 1.019389096537539
$ws.Range("I17").Value = 1.029563697276315
$ws.Range("J17").Value = 1.011448405046783
$ws.Range("K17").Value = 1.030171443724583
$ws.Range("L17").Value = 1.011586414153079
$ws.Range("M17").Value = 1.022813866999529
$ws.Range("N17").Value = 1.007933551255143

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.005274952278349
$ws.Range("D18").Value = 1.026857377938861
$ws.Range("E18").Value = 1.008309597633056
$ws.Range("F18").Value = 1.019655643645727
$ws.Range("I18").Value = 1.029574464784205
$ws.Range("J18").Value = 1.011588389809871
$ws.Range("K18").Value = 1.030213061568839
$ws.Range("L18").Value = 1.011731644282511
$ws.Range("M18").Value = 1.023036778440883
$ws.Range("N18").Value = 1.007981168705495

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.005351171318009
$ws.Range("D19").Value = 1.026886101109538
$ws.Range("E19").Value = 1.008374111786643
$ws.Range("F19").Value = 1.01974653035401
$ws.Range("I19").Value = 1.029578105553583
$ws.Range("J19").Value = 1.011636121799599
$ws.Range("K19").Value = 1.030227233395076
$ws.Range("L19").Value = 1.011781168790274
$ws.Range("M19").Value = 1.023112775347161
$ws.Range("N19").Value = 1.00799740256348

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.005010347283887
$ws.Range("D20").Value = 1.026757630422134
$ws.Range("E20").Value = 1.008085659179791
$ws.Range("F20").Value = 1.019340067689683
$ws.Range("I20").Value = 1.029561702120742
$ws.Range("J20").Value = 1.011422656285945
$ws.Range("K20").Value = 1.030163779529218
$ws.Range("L20").Value = 1.011559702438446
$ws.Range("M20").Value = 1.022772859305866
$ws.Range("N20").Value = 1.007924791216287

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.003903093503159
$ws.Range("D21").Value = 1.026339724229104
$ws.Range("E21").Value = 1.007149101637055
$ws.Range("F21").Value = 1.018018688705476
$ws.Range("I21").Value = 1.029506250640667
$ws.Range("J21").Value = 1.010728716385341
$ws.Range("K21").Value = 1.029956188996813
$ws.Range("L21").Value = 1.010840027721789
$ws.Range("M21").Value = 1.021667057161096
$ws.Range("N21").Value = 1.007688555087603

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.003207531305793
$ws.Range("D22").Value = 1.026076799982181
$ws.Range("E22").Value = 1.006561197304721
$ws.Range("F22").Value = 1.017187937964403
$ws.Range("I22").Value = 1.029469762904783
$ws.Range("J22").Value = 1.010292454363565
$ws.Range("K22").Value = 1.029824681952121
$ws.Range("L22").Value = 1.010387798038532
$ws.Range("M22").Value = 1.020971257231229
$ws.Range("N22").Value = 1.007539895204705

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.003576203927677
$ws.Range("D23").Value = 1.026216196538512
$ws.Range("E23").Value = 1.006872766536364
$ws.Range("F23").Value = 1.017628329743178
$ws.Range("I23").Value = 1.029489259371243
$ws.Range("J23").Value = 1.010523720685779
$ws.Range("K23").Value = 1.029894489416517
$ws.Range("L23").Value = 1.01062750889651
$ws.Range("M23").Value = 1.021340164955155
$ws.Range("N23").Value = 1.007618714759946

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.005028920913016
$ws.Range("D24").Value = 1.026764633641311
$ws.Range("E24").Value = 1.008101376681624
$ws.Range("F24").Value = 1.019362221696668
$ws.Range("I24").Value = 1.029562604204099
$ws.Range("J24").Value = 1.011434291027779
$ws.Range("K24").Value = 1.030167242990805
$ws.Range("L24").Value = 1.011571772225026
$ws.Range("M24").Value = 1.02279138910188
$ws.Range("N24").Value = 1.007928749545506

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.006717380904985
$ws.Range("D25").Value = 1.027400255263292
$ws.Range("E25").Value = 1.009531206342729
$ws.Range("F25").Value = 1.021374544595279
$ws.Range("I25").Value = 1.029640625409049
$ws.Range("J25").Value = 1.012491144060359
$ws.Range("K25").Value = 1.030479418675994
$ws.Range("L25").Value = 1.012668645916531
$ws.Range("M25").Value = 1.024473109964403
$ws.Range("N25").Value = 1.008287959003933

